$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 26.81310566666667
$ws.Range("H2").Value = 80.439317
$ws.Range("I2").Value = 0.004518206005002021
$ws.Range("J2").Value = 0.004518206005002021
$ws.Range("M2").Value = 281.0920463333333
$ws.Range("N2").Value = 843.2761389999999
$ws.Range("O2").Value = 0.8291026083535286
$ws.Range("P2").Value = 0.8291026083535286
$ws.Range("Q2").Value = 7536.950740395229
$ws.Range("R2").Value = 67832.55666355706
$ws.Range("S2").Value = 0.003746056383825751
$ws.Range("T2").Value = 0.003746056383825751

$ws.Range("G3").Value = 26.81310566666667
$ws.Range("H3").Value = 80.439317
$ws.Range("I3").Value = 0.004518206005002021
$ws.Range("J3").Value = 0.004518206005002021
$ws.Range("O3").Value = 0.001324719879221983
$ws.Range("P3").Value = 0.001324719879221983
$ws.Range("Q3").Value = 12.04235564322467
$ws.Range("R3").Value = 108.381200789022
$ws.Range("S3").Value = 0.000005985357313246317
$ws.Range("T3").Value = 0.000005985357313246317

$ws.Range("G4").Value = 26.81310566666667
$ws.Range("H4").Value = 80.439317
$ws.Range("I4").Value = 0.004518206005002021
$ws.Range("J4").Value = 0.004518206005002021
$ws.Range("M4").Value = 4.452417
$ws.Range("N4").Value = 13.357251
$ws.Range("O4").Value = 0.01313274635953239
$ws.Range("P4").Value = 0.01313274635953239
$ws.Range("Q4").Value = 119.383127493063
$ws.Range("R4").Value = 1074.448147437567
$ws.Range("S4").Value = 0.00005933645346380766
$ws.Range("T4").Value = 0.00005933645346380766

$ws.Range("G5").Value = 26.81310566666667
$ws.Range("H5").Value = 80.439317
$ws.Range("I5").Value = 0.004518206005002021
$ws.Range("J5").Value = 0.004518206005002021
$ws.Range("M5").Value = 53.03808999999999
$ws.Range("N5").Value = 159.11427
$ws.Range("O5").Value = 0.156439925407717
$ws.Range("P5").Value = 0.156439925407717
$ws.Range("Q5").Value = 1422.115911528176
$ws.Range("R5").Value = 12799.04320375359
$ws.Range("S5").Value = 0.0007068278103992151
$ws.Range("T5").Value = 0.0007068278103992152

$ws.Range("G6").Value = 5771.873535333333
$ws.Range("I6").Value = 0.9726032482643521
$ws.Range("J6").Value = 0.9726032482643523
$ws.Range("M6").Value = 281.0920463333333
$ws.Range("N6").Value = 843.2761389999999
$ws.Range("O6").Value = 0.8291026083535286
$ws.Range("P6").Value = 0.8291026083535286
$ws.Range("Q6").Value = 1622427.743224058
$ws.Range("R6").Value = 14601849.68901652
$ws.Range("S6").Value = 0.8063878900290888
$ws.Range("T6").Value = 0.8063878900290889

$ws.Range("G7").Value = 5771.873535333333
$ws.Range("I7").Value = 0.9726032482643521
$ws.Range("J7").Value = 0.9726032482643523
$ws.Range("O7").Value = 0.001324719879221983
$ws.Range("P7").Value = 0.001324719879221983
$ws.Range("R7").Value = 23330.4784734238
$ws.Range("S7").Value = 0.001288426857571661
$ws.Range("T7").Value = 0.001288426857571662

$ws.Range("G8").Value = 5771.873535333333
$ws.Range("I8").Value = 0.9726032482643521
$ws.Range("J8").Value = 0.9726032482643523
$ws.Range("M8").Value = 4.452417
$ws.Range("N8").Value = 13.357251
$ws.Range("O8").Value = 0.01313274635953239
$ws.Range("P8").Value = 0.01313274635953239
$ws.Range("Q8").Value = 25698.78785056823
$ws.Range("R8").Value = 231289.0906551141
$ws.Range("S8").Value = 0.01277295176791305
$ws.Range("T8").Value = 0.01277295176791305

$ws.Range("G9").Value = 5771.873535333333
$ws.Range("I9").Value = 0.9726032482643521
$ws.Range("J9").Value = 0.9726032482643523
$ws.Range("M9").Value = 53.03808999999999
$ws.Range("N9").Value = 159.11427
$ws.Range("O9").Value = 0.156439925407717
$ws.Range("P9").Value = 0.156439925407717
$ws.Range("Q9").Value = 306129.1480356274
$ws.Range("R9").Value = 2755162.332320647
$ws.Range("S9").Value = 0.1521539796097785
$ws.Range("T9").Value = 0.1521539796097785

$ws.Range("G10").Value = 132.4457753333333
$ws.Range("H10").Value = 397.337326
$ws.Range("I10").Value = 0.02231808970163987
$ws.Range("J10").Value = 0.02231808970163988
$ws.Range("M10").Value = 281.0920463333333
$ws.Range("N10").Value = 843.2761389999999
$ws.Range("O10").Value = 0.8291026083535286
$ws.Range("P10").Value = 0.8291026083535286
$ws.Range("Q10").Value = 37229.45401665159
$ws.Range("R10").Value = 335065.0861498643
$ws.Range("S10").Value = 0.01850398638509764
$ws.Range("T10").Value = 0.01850398638509764

$ws.Range("G11").Value = 132.4457753333333
$ws.Range("H11").Value = 397.337326
$ws.Range("I11").Value = 0.02231808970163987
$ws.Range("J11").Value = 0.02231808970163988
$ws.Range("O11").Value = 0.001324719879221983
$ws.Range("P11").Value = 0.001324719879221983
$ws.Range("Q11").Value = 59.48431150925733
$ws.Range("R11").Value = 535.3588035833161
$ws.Range("S11").Value = 0.00002956521709402177
$ws.Range("T11").Value = 0.00002956521709402177

$ws.Range("G12").Value = 132.4457753333333
$ws.Range("H12").Value = 397.337326
$ws.Range("I12").Value = 0.02231808970163987
$ws.Range("J12").Value = 0.02231808970163988
$ws.Range("M12").Value = 4.452417
$ws.Range("N12").Value = 13.357251
$ws.Range("O12").Value = 0.01313274635953239
$ws.Range("P12").Value = 0.01313274635953239
$ws.Range("Q12").Value = 589.7038216723139
$ws.Range("R12").Value = 5307.334395050826
$ws.Range("S12").Value = 0.0002930978112809283
$ws.Range("T12").Value = 0.0002930978112809284

$ws.Range("G13").Value = 132.4457753333333
$ws.Range("H13").Value = 397.337326
$ws.Range("I13").Value = 0.02231808970163987
$ws.Range("J13").Value = 0.02231808970163988
$ws.Range("M13").Value = 53.03808999999999
$ws.Range("N13").Value = 159.11427
$ws.Range("O13").Value = 0.156439925407717
$ws.Range("P13").Value = 0.156439925407717
$ws.Range("Q13").Value = 7024.670952249112
$ws.Range("R13").Value = 63222.03857024202
$ws.Range("S13").Value = 0.003491440288167278
$ws.Range("T13").Value = 0.00349144028816728

$ws.Range("G14").Value = 3.326003
$ws.Range("H14").Value = 9.978009
$ws.Range("I14").Value = 0.0005604560290058679
$ws.Range("J14").Value = 0.000560456029005868
$ws.Range("M14").Value = 281.0920463333333
$ws.Range("N14").Value = 843.2761389999999
$ws.Range("O14").Value = 0.8291026083535286
$ws.Range("P14").Value = 0.8291026083535286
$ws.Range("Q14").Value = 934.9129893808056
$ws.Range("R14").Value = 8414.21690442725
$ws.Range("S14").Value = 0.0004646755555162259
$ws.Range("T14").Value = 0.000464675555516226

$ws.Range("G15").Value = 3.326003
$ws.Range("H15").Value = 9.978009
$ws.Range("I15").Value = 0.0005604560290058679
$ws.Range("J15").Value = 0.000560456029005868
$ws.Range("O15").Value = 0.001324719879221983
$ws.Range("P15").Value = 0.001324719879221983
$ws.Range("Q15").Value = 1.493781119366
$ws.Range("R15").Value = 13.444030074294
$ws.Range("S15").Value = 0.0000007424472430538857
$ws.Range("T15").Value = 0.0000007424472430538859

$ws.Range("G16").Value = 3.326003
$ws.Range("H16").Value = 9.978009
$ws.Range("I16").Value = 0.0005604560290058679
$ws.Range("J16").Value = 0.000560456029005868
$ws.Range("M16").Value = 4.452417
$ws.Range("N16").Value = 13.357251
$ws.Range("O16").Value = 0.01313274635953239
$ws.Range("P16").Value = 0.01313274635953239
$ws.Range("Q16").Value = 14.808752299251
$ws.Range("R16").Value = 133.278770693259
$ws.Range("S16").Value = 0.00000736032687460479
$ws.Range("T16").Value = 0.000007360326874604792

$ws.Range("G17").Value = 3.326003
$ws.Range("H17").Value = 9.978009
$ws.Range("I17").Value = 0.0005604560290058679
$ws.Range("J17").Value = 0.000560456029005868
$ws.Range("M17").Value = 53.03808999999999
$ws.Range("N17").Value = 159.11427
$ws.Range("O17").Value = 0.156439925407717
$ws.Range("P17").Value = 0.156439925407717
$ws.Range("Q17").Value = 176.40484645427
$ws.Range("R17").Value = 1587.64361808843
$ws.Range("S17").Value = 0.00008767769937198324
$ws.Range("T17").Value = 0.00008767769937198327
